$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "auswählen"
$ws.Range("B17").Value = "das Selektieren einer dargebotenen Option auf der Nutzeroberfläche"
